$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These metric cells hold text that looks numeric ("0.78", "1.00", ...).
# A leading apostrophe forces Excel to store the literal text (as a
# shared-string cell) instead of auto-converting the entry to a number.
# Only cells whose displayed text actually changes are touched below;
# cells that already show the same text are left alone.

# Row 2
$ws.Range("A2").Value = "'0.78"
$ws.Range("C2").Value = "'0.85"

# Row 3
$ws.Range("B3").Value = "'0.67"
$ws.Range("C3").Value = "'0.80"

# Row 4
$ws.Range("B4").Value = "'0.80"
$ws.Range("C4").Value = "'0.89"

# Row 5
$ws.Range("B5").Value = "'0.73"
$ws.Range("C5").Value = "'0.85"

# Row 6
$ws.Range("B6").Value = "'0.93"
$ws.Range("C6").Value = "'0.97"

# Row 7
$ws.Range("A7").Value = "'0.87"
$ws.Range("B7").Value = "'0.87"
$ws.Range("C7").Value = "'0.87"

# Row 8
$ws.Range("A8").Value = "'0.91"
$ws.Range("B8").Value = "'0.67"
$ws.Range("C8").Value = "'0.77"

# Row 9
$ws.Range("A9").Value = "'0.88"
$ws.Range("B9").Value = "'0.93"
$ws.Range("C9").Value = "'0.90"

# Row 10
$ws.Range("A10").Value = "'0.82"
$ws.Range("B10").Value = "'0.93"
$ws.Range("C10").Value = "'0.87"

# Row 11
$ws.Range("A11").Value = "'0.78"
$ws.Range("B11").Value = "'0.93"
$ws.Range("C11").Value = "'0.85"
